$p = $ppt.ActivePresentation
$d = $p.Designs.Item(1)
Write-Output "Before: $($d.Name)"
try {
  $d.Name = "Default"
  Write-Output "After: $($d.Name)"
} catch {
  Write-Output "err: $_"
}
